$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for rule R10 (row 8, column E)
$ws.Range("E8").Value = "GIT UPDATE"

# Match the author's final selection
$ws.Range("E8").Select()
